# Weekly update: a new price record (row) is reported for Feria Lagunitas de
# Puerto Montt - Alcachofa, and it belongs chronologically before the
# existing row 10, so a new row is inserted at position 10 and the former
# rows 10-14 shift down to 11-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 10, pushing the old rows 10-14 down to 11-15.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new market record.
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C10").Value = "Los Lagos"
$ws.Range("D10").Value = 44488
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 100112013
$ws.Range("G10").Value = "Alcachofa"
$ws.Range("H10").Value = "Madrigal"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 12000
$ws.Range("N10").Value = "`$/caja 40 unidades"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 300
$ws.Range("Q10").Value = 40
$ws.Range("R10").Value = "Hortaliza"
